# Atualizacao de bases das ligas: re-sync betting odds rows 198-205
# (row ids/div/date stay put; match+odds data reshuffled among the rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B198").Value = 6236255
$ws.Range("F198").Value = "Deportivo Rayo Zuliano"
$ws.Range("G198").Value = "Caracas"
$ws.Range("H198").Value = 0
$ws.Range("I198").Value = 0
$ws.Range("J198").Value = "D"
$ws.Range("K198").Value = 3.75
$ws.Range("L198").Value = 3.1
$ws.Range("M198").Value = 1.95
$ws.Range("N198").Value = 2.9
$ws.Range("O198").Value = 2.875
$ws.Range("P198").Value = 2.45
$ws.Range("Q198").Value = 0.25
$ws.Range("R198").Value = 1.775
$ws.Range("S198").Value = 2.025
$ws.Range("T198").Value = 2.25
$ws.Range("U198").Value = 1.85
$ws.Range("V198").Value = 1.95
$ws.Range("W198").Value = -1
$ws.Range("X198").Value = 1.875
$ws.Range("Y198").Value = -1
$ws.Range("Z198").Value = 0.3875
$ws.Range("AA198").Value = -0.5
$ws.Range("AB198").Value = -1
$ws.Range("AC198").Value = 0.95

$ws.Range("B199").Value = 6236251
$ws.Range("F199").Value = "Angostura FC"
$ws.Range("G199").Value = "Portuguesa"
$ws.Range("H199").Value = 1
$ws.Range("I199").Value = 2
$ws.Range("J199").Value = "A"
$ws.Range("K199").Value = 3.1
$ws.Range("L199").Value = 3.2
$ws.Range("M199").Value = 2.15
$ws.Range("N199").Value = 4
$ws.Range("O199").Value = 3.6
$ws.Range("P199").Value = 1.75
$ws.Range("Q199").Value = 0.75
$ws.Range("R199").Value = 1.8
$ws.Range("S199").Value = 2
$ws.Range("T199").Value = 2.5
$ws.Range("U199").Value = 1.95
$ws.Range("V199").Value = 1.85
$ws.Range("W199").Value = -1
$ws.Range("X199").Value = -1
$ws.Range("Y199").Value = 0.75
$ws.Range("Z199").Value = -0.5
$ws.Range("AA199").Value = 0.5
$ws.Range("AB199").Value = 0.95
$ws.Range("AC199").Value = -1

$ws.Range("B200").Value = 6236254
$ws.Range("F200").Value = "Academia Puerto Cabello"
$ws.Range("G200").Value = "Estudiantes Merida"
$ws.Range("H200").Value = 1
$ws.Range("I200").Value = 0
$ws.Range("J200").Value = "H"
$ws.Range("K200").Value = 1.727
$ws.Range("L200").Value = 3.4
$ws.Range("M200").Value = 4.333
$ws.Range("N200").Value = 1.666
$ws.Range("O200").Value = 3.4
$ws.Range("P200").Value = 4.75
$ws.Range("Q200").Value = -0.75
$ws.Range("R200").Value = 1.875
$ws.Range("S200").Value = 1.925
$ws.Range("T200").Value = 2.5
$ws.Range("U200").Value = 1.9
$ws.Range("V200").Value = 1.9
$ws.Range("W200").Value = 0.6659999999999999
$ws.Range("X200").Value = -1
$ws.Range("Y200").Value = -1
$ws.Range("Z200").Value = 0.4375
$ws.Range("AA200").Value = -0.5
$ws.Range("AB200").Value = -1
$ws.Range("AC200").Value = 0.8999999999999999

$ws.Range("B202").Value = 6236252
$ws.Range("F202").Value = "Deportivo Tachira"
$ws.Range("G202").Value = "CD Hermanos Colmenares"
$ws.Range("H202").Value = 1
$ws.Range("I202").Value = 0
$ws.Range("J202").Value = "H"
$ws.Range("K202").Value = 1.363
$ws.Range("L202").Value = 4.2
$ws.Range("M202").Value = 7.5
$ws.Range("N202").Value = 1.333
$ws.Range("O202").Value = 4.5
$ws.Range("P202").Value = 8
$ws.Range("Q202").Value = -1.5
$ws.Range("R202").Value = 2
$ws.Range("S202").Value = 1.8
$ws.Range("T202").Value = 2.5
$ws.Range("U202").Value = 1.925
$ws.Range("V202").Value = 1.875
$ws.Range("W202").Value = 0.333
$ws.Range("X202").Value = -1
$ws.Range("Y202").Value = -1
$ws.Range("Z202").Value = -1
$ws.Range("AA202").Value = 0.8
$ws.Range("AB202").Value = -1
$ws.Range("AC202").Value = 0.875

$ws.Range("B203").Value = 6236612
$ws.Range("F203").Value = "Zamora"
$ws.Range("G203").Value = "Carabobo"
$ws.Range("H203").Value = 0
$ws.Range("I203").Value = 2
$ws.Range("J203").Value = "A"
$ws.Range("K203").Value = 3.2
$ws.Range("L203").Value = 3.1
$ws.Range("M203").Value = 2.15
$ws.Range("N203").Value = 4.5
$ws.Range("O203").Value = 3.3
$ws.Range("P203").Value = 1.75
$ws.Range("Q203").Value = 0.5
$ws.Range("R203").Value = 2
$ws.Range("S203").Value = 1.8
$ws.Range("T203").Value = 2.25
$ws.Range("U203").Value = 1.925
$ws.Range("V203").Value = 1.875
$ws.Range("W203").Value = -1
$ws.Range("X203").Value = -1
$ws.Range("Y203").Value = 0.75
$ws.Range("Z203").Value = -1
$ws.Range("AA203").Value = 0.8
$ws.Range("AB203").Value = -0.5
$ws.Range("AC203").Value = 0.4375

$ws.Range("B204").Value = 6236257
$ws.Range("F204").Value = "CD Hermanos Colmenares"
$ws.Range("G204").Value = "Zamora"
$ws.Range("H204").Value = 0
$ws.Range("I204").Value = 2
$ws.Range("J204").Value = "A"
$ws.Range("K204").Value = 2.3
$ws.Range("L204").Value = 3.2
$ws.Range("M204").Value = 2.8
$ws.Range("N204").Value = 1.666
$ws.Range("O204").Value = 3.8
$ws.Range("P204").Value = 4.2
$ws.Range("Q204").Value = -0.75
$ws.Range("R204").Value = 1.9
$ws.Range("S204").Value = 1.9
$ws.Range("T204").Value = 2.75
$ws.Range("U204").Value = 1.9
$ws.Range("V204").Value = 1.9
$ws.Range("W204").Value = -1
$ws.Range("X204").Value = -1
$ws.Range("Y204").Value = 3.2
$ws.Range("Z204").Value = -1
$ws.Range("AA204").Value = 0.8999999999999999
$ws.Range("AB204").Value = -1
$ws.Range("AC204").Value = 0.8999999999999999

$ws.Range("B205").Value = 6236614
$ws.Range("F205").Value = "Mineros"
$ws.Range("G205").Value = "Angostura FC"
$ws.Range("H205").Value = 1
$ws.Range("I205").Value = 2
$ws.Range("J205").Value = "A"
$ws.Range("K205").Value = 2.45
$ws.Range("L205").Value = 3.3
$ws.Range("M205").Value = 2.55
$ws.Range("N205").Value = 1.8
$ws.Range("O205").Value = 3.75
$ws.Range("P205").Value = 3.6
$ws.Range("Q205").Value = -0.5
$ws.Range("R205").Value = 1.825
$ws.Range("S205").Value = 1.975
$ws.Range("T205").Value = 2.75
$ws.Range("U205").Value = 1.8
$ws.Range("V205").Value = 2
$ws.Range("W205").Value = -1
$ws.Range("X205").Value = -1
$ws.Range("Y205").Value = 2.6
$ws.Range("Z205").Value = -1
$ws.Range("AA205").Value = 0.9750000000000001
$ws.Range("AB205").Value = 0.4
$ws.Range("AC205").Value = -0.5
